$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 12471.429  # H21: 13100 -> 12471.429
$ws.Cells.Item(21, 9).Value = 5000  # I21: 8000 -> 5000
$ws.Cells.Item(21, 10).Value = 13716.667  # J21: 13610 -> 13716.667
$ws.Cells.Item(21, 11).Value = 5000  # K21: 8000 -> 5000
$ws.Cells.Item(21, 12).Value = 13716.667  # L21: 13610 -> 13716.667
$ws.Cells.Item(21, 13).Value = -4532  # M21: -7532 -> -4532
$ws.Cells.Item(21, 14).Value = -14652.667  # N21: -14546 -> -14652.667
$ws.Cells.Item(23, 8).Value = 12471.429  # H23: 13100 -> 12471.429
$ws.Cells.Item(23, 9).Value = 5000  # I23: 8000 -> 5000
$ws.Cells.Item(23, 10).Value = 13716.667  # J23: 13610 -> 13716.667
$ws.Cells.Item(23, 11).Value = 5000  # K23: 8000 -> 5000
$ws.Cells.Item(23, 12).Value = 13716.667  # L23: 13610 -> 13716.667
$ws.Cells.Item(23, 13).Value = -4766  # M23: -7766 -> -4766
$ws.Cells.Item(23, 14).Value = -14184.667  # N23: -14078 -> -14184.667
$ws.Cells.Item(55, 8).Value = 529.5454999999999  # H55: 650.8889 -> 529.5454999999999
$ws.Cells.Item(55, 9).Value = 644  # I55: 786.25 -> 644
$ws.Cells.Item(55, 10).Value = 434.16666  # J55: 542.6 -> 434.16666
$ws.Cells.Item(55, 11).Value = 644  # K55: 786.25 -> 644
$ws.Cells.Item(55, 12).Value = 434.16666  # L55: 542.6 -> 434.16666
$ws.Cells.Item(55, 13).Value = -430  # M55: -572.25 -> -430
$ws.Cells.Item(55, 14).Value = -862.16666  # N55: -970.6 -> -862.16666
$ws.Cells.Item(98, 8).Value = 2818.1  # H98: 2834.6365 -> 2818.1
$ws.Cells.Item(98, 9).Value = 2242.3333  # I98: 2397.625 -> 2242.3333
$ws.Cells.Item(98, 10).Value = 8000  # J98: 4000 -> 8000
$ws.Cells.Item(98, 11).Value = 2242.3333  # K98: 2397.625 -> 2242.3333
$ws.Cells.Item(98, 12).Value = 8000  # L98: 4000 -> 8000
$ws.Cells.Item(98, 13).Value = -744.3332999999998  # M98: -899.625 -> -744.3332999999998
$ws.Cells.Item(98, 14).Value = -10996  # N98: -6996 -> -10996
$ws.Cells.Item(107, 8).Value = 427.45  # H107: 476.47827 -> 427.45
$ws.Cells.Item(107, 9).Value = 397.3158  # I107: 438.41177 -> 397.3158
$ws.Cells.Item(107, 10).Value = 1000  # J107: 584.3333 -> 1000
$ws.Cells.Item(107, 11).Value = 397.3158  # K107: 438.41177 -> 397.3158
$ws.Cells.Item(107, 12).Value = 1000  # L107: 584.3333 -> 1000
$ws.Cells.Item(107, 13).Value = 1522.6842  # M107: 1481.58823 -> 1522.6842
$ws.Cells.Item(107, 14).Value = -4840  # N107: -4424.3333 -> -4840
$ws.Cells.Item(113, 8).Value = 73232.5  # H113: 73282.86 -> 73232.5
$ws.Cells.Item(113, 9).Value = 201471  # I113: 334734.66 -> 201471
$ws.Cells.Item(113, 10).Value = 1988.8889  # J113: 1977.8182 -> 1988.8889
$ws.Cells.Item(113, 11).Value = 201471  # K113: 334734.66 -> 201471
$ws.Cells.Item(113, 12).Value = 1988.8889  # L113: 1977.8182 -> 1988.8889
$ws.Cells.Item(113, 13).Value = -198217  # M113: -331480.66 -> -198217
$ws.Cells.Item(113, 14).Value = -8496.8889  # N113: -8485.8182 -> -8496.8889
$ws.Cells.Item(118, 8).Value = 17731.666  # H118: 17733.334 -> 17731.666
$ws.Cells.Item(118, 9).Value = 34163.332  # I118: 50650 -> 34163.332
$ws.Cells.Item(118, 10).Value = 1300  # J118: 1275 -> 1300
$ws.Cells.Item(118, 11).Value = 102489.996  # K118: 151950 -> 102489.996
$ws.Cells.Item(118, 12).Value = 3900  # L118: 3825 -> 3900
$ws.Cells.Item(118, 13).Value = -100832.996  # M118: -150293 -> -100832.996
$ws.Cells.Item(118, 14).Value = -7214  # N118: -7139 -> -7214
$ws.Cells.Item(122, 8).Value = 2818.1  # H122: 2834.6365 -> 2818.1
$ws.Cells.Item(122, 9).Value = 2242.3333  # I122: 2397.625 -> 2242.3333
$ws.Cells.Item(122, 10).Value = 8000  # J122: 4000 -> 8000
$ws.Cells.Item(122, 11).Value = 6726.999899999999  # K122: 7192.875 -> 6726.999899999999
$ws.Cells.Item(122, 12).Value = 24000  # L122: 12000 -> 24000
$ws.Cells.Item(122, 13).Value = -4276.999899999999  # M122: -4742.875 -> -4276.999899999999
$ws.Cells.Item(122, 14).Value = -28900  # N122: -16900 -> -28900
$ws.Cells.Item(138, 8).Value = 7782.237  # H138: 8057.838 -> 7782.237
$ws.Cells.Item(138, 9).Value = 1440.9  # I138: 1508.6428 -> 1440.9
$ws.Cells.Item(138, 10).Value = 31562.25  # J138: 28433.111 -> 31562.25
$ws.Cells.Item(138, 11).Value = 4322.700000000001  # K138: 4525.928400000001 -> 4322.700000000001
$ws.Cells.Item(138, 12).Value = 94686.75  # L138: 85299.333 -> 94686.75
$ws.Cells.Item(138, 13).Value = 817.2999999999993  # M138: 614.0715999999993 -> 817.2999999999993
$ws.Cells.Item(138, 14).Value = -104966.75  # N138: -95579.333 -> -104966.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 3333.8333  # H12: 3250.75 -> 3333.8333
$ws.Cells.Item(12, 10).Value = 3800  # J12: 4000 -> 3800
$ws.Cells.Item(12, 12).Value = 3800  # L12: 4000 -> 3800
$ws.Cells.Item(12, 14).Value = -4146  # N12: -4346 -> -4146
$ws.Cells.Item(28, 8).Value = 17106.666  # H28: 22044.572 -> 17106.666
$ws.Cells.Item(28, 9).Value = 4413.6665  # I28: 18942.4 -> 4413.6665
$ws.Cells.Item(28, 10).Value = 29799.666  # J28: 29800 -> 29799.666
$ws.Cells.Item(28, 11).Value = 4413.6665  # K28: 18942.4 -> 4413.6665
$ws.Cells.Item(28, 12).Value = 29799.666  # L28: 29800 -> 29799.666
$ws.Cells.Item(28, 13).Value = -4221.6665  # M28: -18750.4 -> -4221.6665
$ws.Cells.Item(28, 14).Value = -30183.666  # N28: -30184 -> -30183.666
$ws.Cells.Item(61, 8).Value = 2581.7273  # H61: 2253.7693 -> 2581.7273
$ws.Cells.Item(61, 9).Value = 2349.75  # I61: 1716.5 -> 2349.75
$ws.Cells.Item(61, 11).Value = 2349.75  # K61: 1716.5 -> 2349.75
$ws.Cells.Item(61, 13).Value = -2137.75  # M61: -1504.5 -> -2137.75
$ws.Cells.Item(74, 8).Value = 736.95123  # H74: 868.8570999999999 -> 736.95123
$ws.Cells.Item(74, 9).Value = 728.4828  # I74: 893 -> 728.4828
$ws.Cells.Item(74, 10).Value = 757.4167  # J74: 796.4286 -> 757.4167
$ws.Cells.Item(74, 11).Value = 728.4828  # K74: 893 -> 728.4828
$ws.Cells.Item(74, 12).Value = 757.4167  # L74: 796.4286 -> 757.4167
$ws.Cells.Item(74, 13).Value = 145.5172  # M74: -19 -> 145.5172
$ws.Cells.Item(74, 14).Value = -2505.4167  # N74: -2544.4286 -> -2505.4167
$ws.Cells.Item(77, 8).Value = 736.95123  # H77: 868.8570999999999 -> 736.95123
$ws.Cells.Item(77, 9).Value = 728.4828  # I77: 893 -> 728.4828
$ws.Cells.Item(77, 10).Value = 757.4167  # J77: 796.4286 -> 757.4167
$ws.Cells.Item(77, 11).Value = 3642.414  # K77: 4465 -> 3642.414
$ws.Cells.Item(77, 12).Value = 3787.0835  # L77: 3982.143 -> 3787.0835
$ws.Cells.Item(77, 13).Value = 725.5860000000002  # M77: -97 -> 725.5860000000002
$ws.Cells.Item(77, 14).Value = -12523.0835  # N77: -12718.143 -> -12523.0835
$ws.Cells.Item(99, 8).Value = 17106.666  # H99: 22044.572 -> 17106.666
$ws.Cells.Item(99, 9).Value = 4413.6665  # I99: 18942.4 -> 4413.6665
$ws.Cells.Item(99, 10).Value = 29799.666  # J99: 29800 -> 29799.666
$ws.Cells.Item(99, 11).Value = 4413.6665  # K99: 18942.4 -> 4413.6665
$ws.Cells.Item(99, 12).Value = 29799.666  # L99: 29800 -> 29799.666
$ws.Cells.Item(99, 13).Value = -1418.6665  # M99: -15947.4 -> -1418.6665
$ws.Cells.Item(99, 14).Value = -35789.666  # N99: -35790 -> -35789.666
$ws.Cells.Item(108, 8).Value = 35000  # H108: 0 -> 35000
$ws.Cells.Item(108, 10).Value = 35000  # J108: 0 -> 35000
$ws.Cells.Item(108, 12).Value = 35000  # L108: 0 -> 35000
$ws.Cells.Item(108, 14).Value = -42680  # N108: None -> -42680
$ws.Cells.Item(122, 8).Value = 1530.1034  # H122: 1524.0358 -> 1530.1034
$ws.Cells.Item(122, 10).Value = 1766.6666  # J122: 1800 -> 1766.6666
$ws.Cells.Item(122, 12).Value = 5299.9998  # L122: 5400 -> 5299.9998
$ws.Cells.Item(122, 14).Value = -10199.9998  # N122: -10300 -> -10199.9998
$ws.Cells.Item(132, 8).Value = 2681.2307  # H132: 2359.6047 -> 2681.2307
$ws.Cells.Item(132, 9).Value = 2455.1428  # I132: 1974.3939 -> 2455.1428
$ws.Cells.Item(132, 11).Value = 7365.428400000001  # K132: 5923.1817 -> 7365.428400000001
$ws.Cells.Item(132, 13).Value = -4835.428400000001  # M132: -3393.1817 -> -4835.428400000001
$ws.Cells.Item(136, 8).Value = 2581.7273  # H136: 2253.7693 -> 2581.7273
$ws.Cells.Item(136, 9).Value = 2349.75  # I136: 1716.5 -> 2349.75
$ws.Cells.Item(136, 11).Value = 7049.25  # K136: 5149.5 -> 7049.25
$ws.Cells.Item(136, 13).Value = -4499.25  # M136: -2599.5 -> -4499.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(42, 8).Value = 398000  # H42: 292000.34 -> 398000
$ws.Cells.Item(42, 10).Value = 398000  # J42: 292000.34 -> 398000
$ws.Cells.Item(42, 12).Value = 398000  # L42: 292000.34 -> 398000
$ws.Cells.Item(42, 14).Value = -398656  # N42: -292656.34 -> -398656
$ws.Cells.Item(134, 8).Value = 3554.102  # H134: 4114.122 -> 3554.102
$ws.Cells.Item(134, 9).Value = 3273.973  # I134: 3988.4482 -> 3273.973
$ws.Cells.Item(134, 11).Value = 9821.919  # K134: 11965.3446 -> 9821.919
$ws.Cells.Item(134, 13).Value = -7286.919  # M134: -9430.3446 -> -7286.919

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 48322.387  # H31: 36974.78 -> 48322.387
$ws.Cells.Item(31, 9).Value = 1500  # I31: 1571.7142 -> 1500
$ws.Cells.Item(31, 10).Value = 51551.516  # J31: 44263.65 -> 51551.516
$ws.Cells.Item(31, 11).Value = 1500  # K31: 1571.7142 -> 1500
$ws.Cells.Item(31, 12).Value = 51551.516  # L31: 44263.65 -> 51551.516
$ws.Cells.Item(31, 13).Value = -1205  # M31: -1276.7142 -> -1205
$ws.Cells.Item(31, 14).Value = -52141.516  # N31: -44853.65 -> -52141.516
$ws.Cells.Item(34, 8).Value = 48322.387  # H34: 36974.78 -> 48322.387
$ws.Cells.Item(34, 9).Value = 1500  # I34: 1571.7142 -> 1500
$ws.Cells.Item(34, 10).Value = 51551.516  # J34: 44263.65 -> 51551.516
$ws.Cells.Item(34, 11).Value = 1500  # K34: 1571.7142 -> 1500
$ws.Cells.Item(34, 12).Value = 51551.516  # L34: 44263.65 -> 51551.516
$ws.Cells.Item(34, 13).Value = -1298  # M34: -1369.7142 -> -1298
$ws.Cells.Item(34, 14).Value = -51955.516  # N34: -44667.65 -> -51955.516
$ws.Cells.Item(58, 8).Value = 1026.8572  # H58: 1096.1333 -> 1026.8572
$ws.Cells.Item(58, 9).Value = 990.97675  # I58: 1048.05 -> 990.97675
$ws.Cells.Item(58, 10).Value = 1284  # J58: 1480.8 -> 1284
$ws.Cells.Item(58, 11).Value = 990.97675  # K58: 1048.05 -> 990.97675
$ws.Cells.Item(58, 12).Value = 1284  # L58: 1480.8 -> 1284
$ws.Cells.Item(58, 13).Value = -787.97675  # M58: -845.05 -> -787.97675
$ws.Cells.Item(58, 14).Value = -1690  # N58: -1886.8 -> -1690
$ws.Cells.Item(132, 8).Value = 30615326  # H132: 33336618 -> 30615326
$ws.Cells.Item(132, 9).Value = 30306268  # I132: 33336818 -> 30306268
$ws.Cells.Item(132, 10).Value = 31252758  # J132: 33336214 -> 31252758
$ws.Cells.Item(132, 11).Value = 90918804  # K132: 100010454 -> 90918804
$ws.Cells.Item(132, 12).Value = 93758274  # L132: 100008642 -> 93758274
$ws.Cells.Item(132, 13).Value = -90916274  # M132: -100007924 -> -90916274
$ws.Cells.Item(132, 14).Value = -93763334  # N132: -100013702 -> -93763334
$ws.Cells.Item(136, 8).Value = 1026.8572  # H136: 1096.1333 -> 1026.8572
$ws.Cells.Item(136, 9).Value = 990.97675  # I136: 1048.05 -> 990.97675
$ws.Cells.Item(136, 10).Value = 1284  # J136: 1480.8 -> 1284
$ws.Cells.Item(136, 11).Value = 2972.93025  # K136: 3144.15 -> 2972.93025
$ws.Cells.Item(136, 12).Value = 3852  # L136: 4442.4 -> 3852
$ws.Cells.Item(136, 13).Value = -422.9302500000003  # M136: -594.1499999999996 -> -422.9302500000003
$ws.Cells.Item(136, 14).Value = -8952  # N136: -9542.4 -> -8952

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 709  # H34: 694.7143 -> 709
$ws.Cells.Item(34, 10).Value = 969.125  # J34: 944.125 -> 969.125
$ws.Cells.Item(34, 12).Value = 2907.375  # L34: 2832.375 -> 2907.375
$ws.Cells.Item(34, 14).Value = -3075.375  # N34: -3000.375 -> -3075.375
$ws.Cells.Item(37, 8).Value = 578143.9  # H37: 1270815 -> 578143.9
$ws.Cells.Item(37, 10).Value = 578143.9  # J37: 1270815 -> 578143.9
$ws.Cells.Item(37, 12).Value = 1734431.7  # L37: 3812445 -> 1734431.7
$ws.Cells.Item(37, 14).Value = -1734655.7  # N37: -3812669 -> -1734655.7
$ws.Cells.Item(97, 8).Value = 932.6667  # H97: 999 -> 932.6667
$ws.Cells.Item(97, 10).Value = 932.6667  # J97: 999 -> 932.6667
$ws.Cells.Item(97, 12).Value = 2798.0001  # L97: 2997 -> 2798.0001
$ws.Cells.Item(97, 14).Value = -3790.0001  # N97: -3989 -> -3790.0001
$ws.Cells.Item(99, 8).Value = 2422.5  # H99: 2450 -> 2422.5
$ws.Cells.Item(99, 9).Value = 1862.5  # I99: 1775 -> 1862.5
$ws.Cells.Item(99, 10).Value = 2609.1667  # J99: 3800 -> 2609.1667
$ws.Cells.Item(99, 11).Value = 5587.5  # K99: 5325 -> 5587.5
$ws.Cells.Item(99, 12).Value = 7827.500100000001  # L99: 11400 -> 7827.500100000001
$ws.Cells.Item(99, 13).Value = -3341.5  # M99: -3079 -> -3341.5
$ws.Cells.Item(99, 14).Value = -12319.5001  # N99: -15892 -> -12319.5001
$ws.Cells.Item(100, 8).Value = 2913.6924  # H100: 2030.6666 -> 2913.6924
$ws.Cells.Item(100, 9).Value = 1180  # I100: 1480 -> 1180
$ws.Cells.Item(100, 10).Value = 3058.1667  # J100: 2080.7273 -> 3058.1667
$ws.Cells.Item(100, 11).Value = 3540  # K100: 4440 -> 3540
$ws.Cells.Item(100, 12).Value = 9174.500100000001  # L100: 6242.1819 -> 9174.500100000001
$ws.Cells.Item(100, 13).Value = -2729  # M100: -3629 -> -2729
$ws.Cells.Item(100, 14).Value = -10796.5001  # N100: -7864.1819 -> -10796.5001
$ws.Cells.Item(101, 8).Value = 3980  # H101: 3975 -> 3980
$ws.Cells.Item(101, 10).Value = 3980  # J101: 3975 -> 3980
$ws.Cells.Item(101, 12).Value = 11940  # L101: 11925 -> 11940
$ws.Cells.Item(101, 14).Value = -16808  # N101: -16793 -> -16808
$ws.Cells.Item(102, 8).Value = 2750  # H102: 4687.6 -> 2750
$ws.Cells.Item(102, 9).Value = 1000  # I102: 4500 -> 1000
$ws.Cells.Item(102, 10).Value = 4500  # J102: 4734.5 -> 4500
$ws.Cells.Item(102, 11).Value = 3000  # K102: 13500 -> 3000
$ws.Cells.Item(102, 12).Value = 13500  # L102: 14203.5 -> 13500
$ws.Cells.Item(102, 13).Value = -566  # M102: -11066 -> -566
$ws.Cells.Item(102, 14).Value = -18368  # N102: -19071.5 -> -18368
$ws.Cells.Item(103, 8).Value = 670.1429000000001  # H103: 940.4545000000001 -> 670.1429000000001
$ws.Cells.Item(103, 9).Value = 670.1429000000001  # I103: 641.6667 -> 670.1429000000001
$ws.Cells.Item(103, 10).Value = 0  # J103: 1299 -> 0
$ws.Cells.Item(103, 11).Value = 2010.4287  # K103: 1925.0001 -> 2010.4287
$ws.Cells.Item(103, 12).Value = 0  # L103: 3897 -> 0
$ws.Cells.Item(103, 13).Value = -1131.4287  # M103: -1046.0001 -> -1131.4287
$ws.Cells.Item(103, 14).ClearContents()  # N103: -5655 -> (removed)
$ws.Cells.Item(104, 8).Value = 3750  # H104: 4500 -> 3750
$ws.Cells.Item(104, 9).Value = 2000  # I104: 0 -> 2000
$ws.Cells.Item(104, 10).Value = 4333.3335  # J104: 4500 -> 4333.3335
$ws.Cells.Item(104, 11).Value = 6000  # K104: 0 -> 6000
$ws.Cells.Item(104, 12).Value = 13000.0005  # L104: 13500 -> 13000.0005
$ws.Cells.Item(104, 13).Value = -3379  # M104: None -> -3379
$ws.Cells.Item(104, 14).Value = -18242.0005  # N104: -18742 -> -18242.0005
$ws.Cells.Item(105, 8).Value = 98794.37  # H105: 7453.8 -> 98794.37
$ws.Cells.Item(105, 10).Value = 98794.37  # J105: 7453.8 -> 98794.37
$ws.Cells.Item(105, 12).Value = 296383.11  # L105: 22361.4 -> 296383.11
$ws.Cells.Item(105, 14).Value = -301625.11  # N105: -27603.4 -> -301625.11
$ws.Cells.Item(106, 8).Value = 2477.7778  # H106: 2437.5 -> 2477.7778
$ws.Cells.Item(106, 10).Value = 2477.7778  # J106: 2437.5 -> 2477.7778
$ws.Cells.Item(106, 12).Value = 7433.3334  # L106: 7312.5 -> 7433.3334
$ws.Cells.Item(106, 14).Value = -9325.3334  # N106: -9204.5 -> -9325.3334
$ws.Cells.Item(131, 8).Value = 6554.99  # H131: 1110.37 -> 6554.99
$ws.Cells.Item(131, 10).Value = 7104.4946  # J131: 1121.3956 -> 7104.4946
$ws.Cells.Item(131, 12).Value = 21313.4838  # L131: 3364.1868 -> 21313.4838
$ws.Cells.Item(131, 14).Value = -31393.4838  # N131: -13444.1868 -> -31393.4838

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 104895.5  # H70: 65146.305 -> 104895.5
$ws.Cells.Item(70, 9).Value = 225990.67  # I70: 128865.25 -> 225990.67
$ws.Cells.Item(70, 10).Value = 5817.636  # J70: 5175.5293 -> 5817.636
$ws.Cells.Item(70, 11).Value = 225990.67  # K70: 128865.25 -> 225990.67
$ws.Cells.Item(70, 12).Value = 5817.636  # L70: 5175.5293 -> 5817.636
$ws.Cells.Item(70, 13).Value = -225720.67  # M70: -128595.25 -> -225720.67
$ws.Cells.Item(70, 14).Value = -6357.636  # N70: -5715.5293 -> -6357.636
$ws.Cells.Item(73, 8).Value = 104895.5  # H73: 65146.305 -> 104895.5
$ws.Cells.Item(73, 9).Value = 225990.67  # I73: 128865.25 -> 225990.67
$ws.Cells.Item(73, 10).Value = 5817.636  # J73: 5175.5293 -> 5817.636
$ws.Cells.Item(73, 11).Value = 225990.67  # K73: 128865.25 -> 225990.67
$ws.Cells.Item(73, 12).Value = 5817.636  # L73: 5175.5293 -> 5817.636
$ws.Cells.Item(73, 13).Value = -225054.67  # M73: -127929.25 -> -225054.67
$ws.Cells.Item(73, 14).Value = -7689.636  # N73: -7047.5293 -> -7689.636
$ws.Cells.Item(97, 8).Value = 83335140  # H97: 66668550 -> 83335140
$ws.Cells.Item(97, 9).Value = 142859580  # I97: 83335370 -> 142859580
$ws.Cells.Item(97, 10).Value = 914.2  # J97: 1270 -> 914.2
$ws.Cells.Item(97, 11).Value = 142859580  # K97: 83335370 -> 142859580
$ws.Cells.Item(97, 12).Value = 914.2  # L97: 1270 -> 914.2
$ws.Cells.Item(97, 13).Value = -142859084  # M97: -83334874 -> -142859084
$ws.Cells.Item(97, 14).Value = -1906.2  # N97: -2262 -> -1906.2
$ws.Cells.Item(102, 8).Value = 4573.3335  # H102: 5193.3335 -> 4573.3335
$ws.Cells.Item(102, 9).Value = 3576  # I102: 3940 -> 3576
$ws.Cells.Item(102, 11).Value = 3576  # K102: 3940 -> 3576
$ws.Cells.Item(102, 13).Value = -1954  # M102: -2318 -> -1954
$ws.Cells.Item(132, 8).Value = 2320.2683  # H132: 2387.8462 -> 2320.2683
$ws.Cells.Item(132, 9).Value = 1614.8276  # I132: 1660.1852 -> 1614.8276
$ws.Cells.Item(132, 11).Value = 4844.4828  # K132: 4980.5556 -> 4844.4828
$ws.Cells.Item(132, 13).Value = -2314.4828  # M132: -2450.5556 -> -2314.4828

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 2500500  # H5: 2524000 -> 2500500
$ws.Cells.Item(5, 10).Value = 2500500  # J5: 2524000 -> 2500500
$ws.Cells.Item(5, 12).Value = 2500500  # L5: 2524000 -> 2500500
$ws.Cells.Item(5, 14).Value = -2500724  # N5: -2524224 -> -2500724
$ws.Cells.Item(18, 8).Value = 0  # H18: 18000 -> 0
$ws.Cells.Item(18, 10).Value = 0  # J18: 18000 -> 0
$ws.Cells.Item(18, 12).Value = 0  # L18: 18000 -> 0
$ws.Cells.Item(18, 14).ClearContents()  # N18: -18346 -> (removed)
$ws.Cells.Item(132, 8).Value = 2404.653  # H132: 2627.568 -> 2404.653
$ws.Cells.Item(132, 9).Value = 2786.353  # I132: 3370.926 -> 2786.353
$ws.Cells.Item(132, 10).Value = 1539.4667  # J132: 1446.9412 -> 1539.4667
$ws.Cells.Item(132, 11).Value = 8359.059000000001  # K132: 10112.778 -> 8359.059000000001
$ws.Cells.Item(132, 12).Value = 4618.4001  # L132: 4340.8236 -> 4618.4001
$ws.Cells.Item(132, 13).Value = -5829.059000000001  # M132: -7582.778 -> -5829.059000000001
$ws.Cells.Item(132, 14).Value = -9678.400099999999  # N132: -9400.8236 -> -9678.400099999999
$ws.Cells.Item(136, 8).Value = 1025.125  # H136: 1112.9615 -> 1025.125
$ws.Cells.Item(136, 9).Value = 604.34784  # I136: 640.6667 -> 604.34784
$ws.Cells.Item(136, 10).Value = 2100.4443  # J136: 2175.625 -> 2100.4443
$ws.Cells.Item(136, 11).Value = 1813.04352  # K136: 1922.0001 -> 1813.04352
$ws.Cells.Item(136, 12).Value = 6301.3329  # L136: 6526.875 -> 6301.3329
$ws.Cells.Item(136, 13).Value = 736.9564799999998  # M136: 627.9999 -> 736.9564799999998
$ws.Cells.Item(136, 14).Value = -11401.3329  # N136: -11626.875 -> -11401.3329
